$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 97 with a new time-log entry (previously a blank placeholder row)
$ws.Range("A97").Value = 41933
$ws.Range("B97").Value = 0.72569444444444453
$ws.Range("C97").Value = 0.76666666666666661
$ws.Range("D97").Value = 5
$ws.Range("F97").Value = "Coding"
$ws.Range("E97").Formula = '=IF(AND(NOT(ISBLANK(B97)),NOT(ISBLANK(C97))), (C97-B97) * 24 - D97/60, "")'

# Move the active selection to reflect where the author last clicked
$ws.Range("C98").Select() | Out-Null
